# Fruta / hortaliza, semanal
# The data rows (2-16) get their Fecha (D), Volumen (J), Precio minimo (K),
# Precio maximo (L), Precio promedio ponderado (M) and Precio $/Kg (P)
# columns reshuffled across rows according to a fixed permutation (the
# weekly re-sequencing of the already-collected daily records). All other
# columns (A, B, C, E, F, G, H, I, N, O, Q, R) stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> row whose old values should be copied into it
$map = @{
    2  = 5
    3  = 6
    4  = 9
    5  = 8
    6  = 16
    7  = 13
    8  = 12
    9  = 3
    10 = 7
    11 = 4
    12 = 15
    13 = 14
    14 = 10
    15 = 2
    16 = 11
}

$cols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P

# Snapshot the current ("before") values for every affected cell first,
# since the permutation reads from multiple rows while writing to others.
$snapshot = @{}
for ($r = 2; $r -le 16; $r++) {
    foreach ($c in $cols) {
        $snapshot["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

foreach ($r in $map.Keys) {
    $src = $map[$r]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $snapshot["$src,$c"]
    }
}
